$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") -- copy H1's formatting (bold,
# centered, bordered header style) onto the new header cells first.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns I and J for rows 2-22.
$values = @{
    2  = @(4, 5)
    3  = @(6, 6)
    4  = @(8, 8)
    5  = @(8, 8)
    6  = @(7, 7)
    7  = @(7, 8)
    8  = @(7, 7)
    9  = @(6, 7)
    10 = @(6, 8)
    11 = @(7, 8)
    12 = @(7, 8)
    13 = @(7, 7)
    14 = @(7, 7)
    15 = @(8, 8)
    16 = @(7, 8)
    17 = @(8, 8)
    18 = @(7, 7)
    19 = @(6, 6)
    20 = @(9, 9)
    21 = @(8, 8)
    22 = @(6, 6)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Range("I$row").Value = $pair[0]
    $ws.Range("J$row").Value = $pair[1]
}
